# Update the EDCR Results worksheet with the new copper results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDCR Results")

# Column H (PRC) changes from 12 to 9 for every data row (2 through 100).
$hRange = $ws.Range("H2:H100")
$hRange.Value = 9

# Column G (NRC) changes from 50 to 40 for rows 18 through 100.
$gRange = $ws.Range("G18:G100")
$gRange.Value = 40
